$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.834.00"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.036.25"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "227.36"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.607"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.45%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "60.10"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.16%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.386"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0817"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  -0.06%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.64"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "2.339.34"
$ws.Range("E13").Value = "  -1.08%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "21.02"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.57%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.761"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.63%  "
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "2.050.37"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "37.775.27"
$ws.Range("E18").Value = "  -0.67%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.07"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.74%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "69.80"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.08%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -0.92%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "224.97"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.10%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.44"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "165.23"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.21"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.32%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.130"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -3.44%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "18.92"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("E30").Value = "  -6.28%  "
$ws.Range("E31").Value = "  +1.19%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.44"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.56%  "
$ws.Range("E33").Value = "  +3.89%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.49"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -2.74%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0601"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.37%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.42"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +5.94%  "
$ws.Range("E37").Value = "  -5.41%  "
$ws.Range("E38").Value = "  -2.57%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "1.539.39"
$ws.Range("E41").Value = "  -0.31%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "97.04"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -1.46%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "16.89"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  -0.86%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0924"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("E46").Value = "  -1.16%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "3.93"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.74%  "
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "2.230.38"
$ws.Range("E51").Value = "  -0.99%  "
